$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $origStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = $origStyle
}

Set-TextValue 'D2' '23.972.13'
$ws.Range('E2').Value = '  -0.49%  '
Set-TextValue 'D3' '1.654.73'
$ws.Range('E3').Value = '  +0.68%  '
Set-TextValue 'D4' '1.003'
$ws.Range('E4').Value = '  -0.10%  '
Set-TextValue 'D5' '309.40'
$ws.Range('E5').Value = '  +0.20%  '
Set-TextValue 'D6' '1.001'
$ws.Range('E6').Value = '  -0.12%  '
Set-TextValue 'D7' '0.3912'
$ws.Range('E7').Value = '  -0.78%  '
Set-TextValue 'D8' '0.3884'
$ws.Range('E8').Value = '  +0.53%  '
Set-TextValue 'D9' '51.84'
$ws.Range('E9').Value = '  +4.49%  '
Set-TextValue 'D10' '1.370'
$ws.Range('E10').Value = '  +0.37%  '
Set-TextValue 'D11' '0.9980'
$ws.Range('E11').Value = '  -0.61%  '
Set-TextValue 'D12' '0.08489'
$ws.Range('E12').Value = '  -0.87%  '
Set-TextValue 'D13' '24.16'
$ws.Range('E13').Value = '  +2.53%  '
Set-TextValue 'D14' '7.264'
$ws.Range('E14').Value = '  +2.57%  '
Set-TextValue 'D15' '8.105'
$ws.Range('E15').Value = '  +7.96%  '
Set-TextValue 'D16' '0.00001317'
$ws.Range('E16').Value = '  +2.57%  '
Set-TextValue 'D17' '1.654.74'
$ws.Range('E17').Value = '  +0.23%  '
Set-TextValue 'D18' '95.17'
$ws.Range('E18').Value = '  +1.37%  '
Set-TextValue 'D19' '0.06965'
$ws.Range('E19').Value = '  +0.78%  '
Set-TextValue 'D20' '19.98'
$ws.Range('E20').Value = '  -1.56%  '
Set-TextValue 'D21' '6.986'
$ws.Range('E21').Value = '  +0.83%  '
Set-TextValue 'D22' '1.002'
$ws.Range('E22').Value = '  +0.08%  '
Set-TextValue 'D23' '13.73'
$ws.Range('E23').Value = '  +0.99%  '
Set-TextValue 'D24' '23.985.68'
$ws.Range('E24').Value = '  -0.51%  '
Set-TextValue 'D25' '3.152'
$ws.Range('E25').Value = '  +9.80%  '
Set-TextValue 'D26' '2.501'
$ws.Range('E26').Value = '  +3.35%  '
Set-TextValue 'D27' '22.31'
$ws.Range('E27').Value = '  +0.62%  '
Set-TextValue 'D28' '153.83'
$ws.Range('E28').Value = '  -2.63%  '
Set-TextValue 'D29' '140.11'
$ws.Range('E29').Value = '  -0.04%  '
Set-TextValue 'D30' '5.292'
$ws.Range('E30').Value = '  +0.74%  '
Set-TextValue 'D31' '7.871'
$ws.Range('E31').Value = '  -3.30%  '
Set-TextValue 'D32' '2.478'
$ws.Range('E32').Value = '  -1.67%  '
Set-TextValue 'D33' '1.836.47'
$ws.Range('E33').Value = '  +0.21%  '
$ws.Range('E34').Value = '  +8.27%  '
Set-TextValue 'D35' '0.03014'
$ws.Range('E35').Value = '  +3.53%  '
Set-TextValue 'D36' '0.08119'
$ws.Range('E36').Value = '  +0.31%  '
Set-TextValue 'D37' '11.20'
$ws.Range('E37').Value = '  +8.04%  '
Set-TextValue 'D38' '6.691'
$ws.Range('E38').Value = '  -1.01%  '
Set-TextValue 'D39' '0.2705'
$ws.Range('E39').Value = '  +0.56%  '
Set-TextValue 'D40' '0.09150'
$ws.Range('E40').Value = '  -0.79%  '
Set-TextValue 'D41' '0.7596'
$ws.Range('E41').Value = '  +1.26%  '
Set-TextValue 'D42' '13.51'
$ws.Range('E42').Value = '  +3.62%  '
$ws.Range('E43').Value = '  -1.19%  '
$ws.Range('E44').Value = '  +2.26%  '
Set-TextValue 'D45' '0.7035'
$ws.Range('E45').Value = '  +2.11%  '
Set-TextValue 'D46' '2.503'
$ws.Range('E46').Value = '  +1.82%  '
Set-TextValue 'D47' '4.083'
$ws.Range('E47').Value = '  -0.29%  '
Set-TextValue 'D48' '1.001'
$ws.Range('E48').Value = '  +0.02%  '
Set-TextValue 'D49' '0.08346'
$ws.Range('E49').Value = '  -0.24%  '
Set-TextValue 'D50' '135.11'
$ws.Range('E50').Value = '  +0.90%  '
Set-TextValue 'D51' '1.239'
$ws.Range('E51').Value = '  -1.94%  '
